$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A - everything that used to live in
# columns A-H now lives in columns B-I.
$ws.Columns("A:A").Insert()

# New column A holds the English field-name labels that correspond to the
# Chinese labels now sitting in column B. Write them in the same order the
# original commit's shared strings appear in (Post, title, content,
# attachment, user, remarks, records) so the generated sharedStrings table
# lines up with the target.
$ws.Range("B21").Value = "Post"
$ws.Range("A4").Value = "title"
$ws.Range("A5").Value = "content"
$ws.Range("A6").Value = "attachment"
$ws.Range("A3").Value = "user"
$ws.Range("A7").Value = "remarks"
$ws.Range("A10").Value = "records"

# Give column A a bit of breathing room for the new labels.
$ws.Columns("A:A").ColumnWidth = 13.5

# Restore the active selection to A10 (post-edit cursor position).
[void]$ws.Range("A10").Select()
